$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Copy the existing date-format (numFmtId 14) style into the brand new
#    column J (rows 2-21) *before* columns D/F lose that style below.
#    The cells stay empty - only the formatting is copied.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("D$r").Copy()
    [void]$ws.Range("J$r").PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 2. promotionDatedate (column F) - turn the date serials into plain text
#    dates formatted as dd/mm/yyyy. Column F first, so the new shared
#    strings are appended in the same order as in the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("F1:F21").NumberFormat = "@"

$ws.Range("F2").Value = "30/09/1996"
$ws.Range("F3").Value = "12/08/1997"
$ws.Range("F4").Value = "30/09/1997"
$ws.Range("F5").Value = "12/08/1998"
$ws.Range("F6").Value = "30/09/1998"
$ws.Range("F7").Value = "12/08/1999"
$ws.Range("F8").Value = "30/09/1999"
$ws.Range("F9").Value = "12/08/2000"
$ws.Range("F10").Value = "30/09/2000"
$ws.Range("F11").Value = "12/08/2001"
$ws.Range("F12").Value = "30/09/2001"
$ws.Range("F13").Value = "12/08/2002"
$ws.Range("F14").Value = "30/09/2002"
$ws.Range("F15").Value = "12/08/2003"
$ws.Range("F16").Value = "30/09/2003"
$ws.Range("F17").Value = "12/08/2004"
$ws.Range("F18").Value = "30/09/2004"
$ws.Range("F19").Value = "12/08/2005"
$ws.Range("F20").Value = "30/09/2005"
$ws.Range("F21").Value = "12/08/2006"

# ---------------------------------------------------------------------------
# 3. joiningDatedate (column D) - same treatment.
# ---------------------------------------------------------------------------
$ws.Range("D1:D21").NumberFormat = "@"

$ws.Range("D2").Value = "12/03/1995"
$ws.Range("D3").Value = "21/12/1995"
$ws.Range("D4").Value = "12/03/1996"
$ws.Range("D5").Value = "21/12/1996"
$ws.Range("D6").Value = "12/03/1997"
$ws.Range("D7").Value = "21/12/1997"
$ws.Range("D8").Value = "12/03/1998"
$ws.Range("D9").Value = "21/12/1998"
$ws.Range("D10").Value = "12/03/1999"
$ws.Range("D11").Value = "21/12/1999"
$ws.Range("D12").Value = "12/03/2000"
$ws.Range("D13").Value = "21/12/2000"
$ws.Range("D14").Value = "12/03/2001"
$ws.Range("D15").Value = "21/12/2001"
$ws.Range("D16").Value = "12/03/2002"
$ws.Range("D17").Value = "21/12/2002"
$ws.Range("D18").Value = "12/03/2003"
$ws.Range("D19").Value = "21/12/2003"
$ws.Range("D20").Value = "12/03/2004"
$ws.Range("D21").Value = "21/12/2004"

# ---------------------------------------------------------------------------
# 4. New column J needs a width, and the workbook's dimension / active
#    selection grow to include it too.
# ---------------------------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 9.6

[void]$ws.Range("H21").Select()
